$wb = $excel.ActiveWorkbook

# Sheet "建物" (Building) - fix property_category column (I) for data rows 2-9
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I9").Value = "building"

# Sheet "汽車" (Car) - fix property_category column (H) for data rows 2-3
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2:H3").Value = "car"
